$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Move Robot40 to location (4, 8) and remove the toolkit."
$ws.Range("B1").Value = "['Robot2']"
$ws.Range("E1").Value = "(4, 8)"

$ws.Range("A2").Value = "Move Robot40 to location (6, 2) and remove the liquid spill."
$ws.Range("B2").Value = "['Robot26']"
$ws.Range("E2").Value = "(6, 2)"

$ws.Range("A3").Value = "Move Robot9 to location (12, 3) and remove the large debris."
$ws.Range("E3").Value = "(12, 3)"

$ws.Range("A4").Value = "Move Robot35 to location (2, 11) and remove the dust."
$ws.Range("B4").Value = "['Robot50', 'Robot28']"
$ws.Range("E4").Value = "(2, 11)"

$ws.Range("A5").Value = "Move Robot26 to location (12, 1) and remove the grass."
$ws.Range("B5").Value = "['Robot31']"
$ws.Range("E5").Value = "(12, 1)"

$ws.Range("A6").Value = "Move Robot41 to location (4, 11) and remove the small debris."
$ws.Range("B6").Value = "['Robot28', 'Robot50']"
$ws.Range("E6").Value = "(4, 11)"

$ws.Range("A7").Value = "Move Robot2 to location (9, 1) and remove the vehicle."
$ws.Range("B7").Value = "['Robot23']"
$ws.Range("E7").Value = "(9, 1)"

$ws.Range("A8").Value = "Move Robot28 to location (11, 6) and remove the construction materials."
$ws.Range("B8").Value = "['Robot42', 'Robot29', 'Robot23']"
$ws.Range("E8").Value = "(11, 6)"

$ws.Range("A9").Value = "Move Robot39 to location (6, 1) and remove the tree branches."
$ws.Range("E9").Value = "(6, 1)"

$ws.Range("A10").Value = "Move Robot28 to location (2, 8) and remove the screws."
$ws.Range("B10").Value = "['Robot45']"
$ws.Range("E10").Value = "(2, 8)"
